# "CAJAS PARA ENRROLLADOR" price list (Hoja1) update:
#  - A1 header date rolls forward one month (24/04/2024 -> 24/05/2024,
#    serial 45406 -> 45436)
#  - D30 price (CAJA p/ ENROLLADOR CHICA) 338 -> 799
#  - D31 price (CAJA p/ ENROLLADOR GRANDE) 405 -> 967

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45436
$ws.Range("D30").Value = 799
$ws.Range("D31").Value = 967
